# Updates the "cryptos" price/volume snapshot on Sheet1 to the values
# recorded by the Wed Oct 11 07:28:43 UTC 2023 GitHub Actions run.
# Column D ("Price") cells whose new text happens to look like a plain
# number are forced to Text first (NumberFormat "@") so Excel doesn't
# silently coerce them into numeric values and drop significant trailing
# zeros (e.g. "206.52" -> 206.5, "0.0460" -> 0.046). ClearFormats()
# afterwards drops the temporary text format again so the cell keeps its
# original (default/general) style.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '27.125.83'
$ws.Range("D3").Value = '1.559.18'
$ws.Range("E3").Value = '  -2.31%  '
$ws.Range("E4").Value = '  -0.01%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '206.52'
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = '  -0.98%  '
$ws.Range("E6").Value = '  -2.90%  '
$ws.Range("E7").Value = '  +0.00%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '22.09'
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = '  -1.36%  '
$ws.Range("E9").Value = '  -2.43%  '
$ws.Range("E10").Value = '  -0.59%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0864'
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = '  -0.49%  '
$ws.Range("D12").Value = '1.781.14'
$ws.Range("E12").Value = '  -2.25%  '
$ws.Range("D13").Value = '1.565.28'
$ws.Range("E13").Value = '  -1.93%  '
$ws.Range("E14").Value = '  -2.89%  '
$ws.Range("E15").Value = '  -3.69%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '62.90'
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = '  -1.05%  '
$ws.Range("D17").Value = '27.132.62'
$ws.Range("E17").Value = '  -2.15%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '215.09'
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = '  -2.30%  '
$ws.Range("E19").Value = '  -1.87%  '
$ws.Range("E20").Value = '  -2.34%  '
$ws.Range("E21").Value = '  +0.01%  '
$ws.Range("E22").Value = '  -1.28%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '9.30'
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = '  -4.57%  '
$ws.Range("E24").Value = '  +0.17%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '151.52'
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = '  -1.63%  '
$ws.Range("E26").Value = '  -3.10%  '
$ws.Range("E27").Value = '  -1.73%  '
$ws.Range("E28").Value = '  +0.04%  '
$ws.Range("E29").Value = '  -1.51%  '
$ws.Range("E30").Value = '  -1.94%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.0460'
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = '  -2.54%  '
$ws.Range("E32").Value = '  -2.54%  '
$ws.Range("D33").Value = '1.382.81'
$ws.Range("E33").Value = '  +0.34%  '
$ws.Range("E34").Value = '  -1.38%  '
$ws.Range("E35").Value = '  -0.60%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.947'
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = '  -2.68%  '
$ws.Range("E37").Value = '  -1.70%  '
$ws.Range("E38").Value = '  -1.91%  '
$ws.Range("E39").Value = '  -2.49%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.514'
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = '  -4.51%  '
$ws.Range("E41").Value = '  +0.03%  '
$ws.Range("E42").Value = '  +1.48%  '
$ws.Range("E43").Value = '  +2.70%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '63.31'
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = '  -2.13%  '
$ws.Range("E45").Value = '  -0.43%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '5.23'
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = '  +0.14%  '
$ws.Range("D47").Value = '1.693.34'
$ws.Range("E47").Value = '  -2.23%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '85.36'
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = '  -1.84%  '
$ws.Range("E49").Value = '  -2.31%  '
$ws.Range("E50").Value = '  -1.03%  '
$ws.Range("E51").Value = '  +0.06%  '
